# Loan RBI, Variable Instalments
# On the "Repayment schedule" sheet, insert a new (blank) column before
# column N so the existing "Late"/"Outstanding"/"Disbursement" columns
# shift one place to the right, give the newly inserted column the same
# width the author gave it, and make "Repayment schedule" the active
# sheet/tab with the selection parked on K16 (as it was left after the
# edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (shifts old N/O/P -> O/P/Q).
$ws.Columns("N").Insert()

# The inserted column keeps a custom width of 11 (no bestFit), matching
# the author's manual column-width tweak after inserting it.
$ws.Columns("N").ColumnWidth = 10.166666666666666

# Make "Repayment schedule" the active sheet/tab and park the selection
# on K16, matching the author's final selection.
$ws.Activate()
$ws.Range("K16").Select()
